$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q4").Value = 1.92
$ws.Range("R4").Value = 1.98
$ws.Range("H5").Value = 3.75
$ws.Range("I5").Value = 2.1
$ws.Range("J5").Value = 3.6
$ws.Range("L5").Value = 2.75
$ws.Range("O5").Value = 1.2
$ws.Range("P5").Value = 3.65
$ws.Range("U5").Value = 1.5
$ws.Range("V5").Value = 2.25
$ws.Range("W5").Value = 13
$ws.Range("Y5").Value = 11
$ws.Range("AG5").Value = 126
$ws.Range("AH5").Value = 10
$ws.Range("AI5").Value = 12
$ws.Range("AO5").Value = 15
$ws.Range("AP5").Value = 21
$ws.Range("AR5").Value = 51
$ws.Range("AT5").Value = 3.4
$ws.Range("AU5").Value = 7
$ws.Range("AX5").Value = 4.5
$ws.Range("AY5").Value = 12
$ws.Range("BA5").Value = 41
$ws.Range("G6").Value = 1.85
$ws.Range("H6").Value = 3.4
$ws.Range("I6").Value = 3.75
$ws.Range("J6").Value = 2.6
$ws.Range("L6").Value = 4.5
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 11
$ws.Range("Q6").Value = 1.98
$ws.Range("R6").Value = 1.83
$ws.Range("U6").Value = 1.83
$ws.Range("V6").Value = 1.83
$ws.Range("W6").Value = 7
$ws.Range("X6").Value = 9
$ws.Range("Z6").Value = 15
$ws.Range("AA6").Value = 15
$ws.Range("AC6").Value = 10
$ws.Range("AI6").Value = 21
$ws.Range("AL6").Value = 34
$ws.Range("AO6").Value = 10
$ws.Range("AQ6").Value = 34
$ws.Range("AX6").Value = 6
$ws.Range("BA6").Value = 81
$ws.Range("BB6").Value = 101
$ws.Range("G7").Value = 6.5
$ws.Range("H7").Value = 5
$ws.Range("I7").Value = 1.33
$ws.Range("J7").Value = 6.5
$ws.Range("L7").Value = 1.8
$ws.Range("M7").Value = 1.01
$ws.Range("N7").Value = 23
$ws.Range("O7").Value = 1.11
$ws.Range("P7").Value = 6.5
$ws.Range("Q7").Value = 1.36
$ws.Range("R7").Value = 3.1
$ws.Range("W7").Value = 26
$ws.Range("Z7").Value = 81
$ws.Range("AC7").Value = 23
$ws.Range("AD7").Value = 11
$ws.Range("AE7").Value = 17
$ws.Range("AK7").Value = 10
$ws.Range("AL7").Value = 10
$ws.Range("AN7").Value = 9
$ws.Range("AO7").Value = 34
$ws.Range("AZ7").Value = 13
$ws.Range("BA7").Value = 15
$ws.Range("BB7").Value = 29
$ws.Range("G14").Value = 25
$ws.Range("H14").Value = 10.5
$ws.Range("I14").Value = 1.04
$ws.Range("J14").Value = 18
$ws.Range("K14").Value = 4.4
$ws.Range("T14").Value = 5.7
$ws.Range("U14").Value = 2.37
$ws.Range("W14").Value = 175
$ws.Range("X14").Value = 700
$ws.Range("Y14").Value = 175
$ws.Range("AA14").Value = 900
$ws.Range("AB14").Value = 450
$ws.Range("AC14").Value = 40
$ws.Range("AF14").Value = 300
$ws.Range("AH14").Value = 21
$ws.Range("AI14").Value = 10.25
$ws.Range("AJ14").Value = 20
$ws.Range("AK14").Value = 7.4
$ws.Range("AL14").Value = 14.5
$ws.Range("AM14").Value = 55
$ws.Range("AN14").Value = 28
$ws.Range("AO14").Value = 200
$ws.Range("AP14").Value = 120
$ws.Range("AT14").Value = 5.7
$ws.Range("AU14").Value = 14.5
$ws.Range("AV14").Value = 110
$ws.Range("AX14").Value = 3.75
$ws.Range("AZ14").Value = 15
$ws.Range("BA14").Value = 5.7
$ws.Range("BC14").Value = 175
$ws.Range("G15").Value = 4.2
$ws.Range("I15").Value = 1.85
$ws.Range("L15").Value = 2.5
$ws.Range("AG15").Value = 301
$ws.Range("AI15").Value = 8.5
$ws.Range("AJ15").Value = 8.5
$ws.Range("AO15").Value = 23
$ws.Range("AS15").Value = 251
$ws.Range("AX15").Value = 3.75
